# changed : make demo data more realistic and complete
# Replace the placeholder "dataset_3" values in the demo table with a more
# realistic dataset name "dep_sante", and move the active selection to B5
# (as reflected in the sheet's saved view state).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "dep_sante"
$ws.Range("B3").Value = "dep_sante"

$ws.Range("B5").Select() | Out-Null
